$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Two new planning rows at the bottom (minted first so the shared-string
#     table ends up in the same order as the canonical export) ---
$ws.Range("A17").Value = "Let characters follow selected character ( no clipping -> navmesh )"
$ws.Range("C6").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "/"

$ws.Range("A18").Value = "Let enemies track ( follow ) selected character ( no clipping -> navmesh ). When close enough, attack closest of 3 characters"
$ws.Range("C6").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "/"

# --- Row 6: clarify "Writing character" task text ---
$ws.Range("A6").Value = "Writing character ( = visualize the prefab and spawn it in level with correct stats )"

# --- New note in D5 (next to the BUSY "Writing managers" row) ---
$ws.Range("D5").Value = "UIManager not entirely clear, UI seems to work without it"

# --- Rows that moved from "/" (not started) to "DONE" (with time spent) ---
# Row 8: Spawn 3 characters -> DONE, 1:00 spent
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "DONE"
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = 0.041666666666666664

# Row 12: Weapon stats altering -> DONE, 0:20 spent
$ws.Range("C2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "DONE"
$ws.Range("B2").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 0.013888888888888888

# Row 13: UI -> DONE, 1:50 spent
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = "DONE"
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 0.0763888888888889

# Row 15: Panels with info ( weapons, health, name ) -> DONE, 0:30 spent
$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "DONE"
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 0.020833333333333332

# --- Widen column A to fit the longer text ---
$ws.Columns.Item(1).ColumnWidth = 110.42

# --- Selection state ends on D5 ---
$ws.Range("D5").Select()
